$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")
$ws.Activate()

# --- Bring formatting for the two new data rows (20 & 21) in line with the
# --- existing "filled" rows they now mirror, before touching values/heights.
$ws.Range("C19:F19").Copy()
$ws.Range("C20:F20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C17:F17").Copy()
$ws.Range("C21:F21").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 20 & 21 become taller (wrapped, two-line task descriptions) -> ht=30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 30

# Column C (Application) now wraps text, same as column D already does.
$ws.Columns.Item(3).WrapText = $true

# --- Re-labelled "Application" tags for several existing rows ---
$ws.Range("C5").Value = "B2C & Hayaai/Sonia"
$ws.Range("C6").Value = "B2C & Hayaai/Sonia"
$ws.Range("C7").Value = "B2C & Hayaai/Sonia"
$ws.Range("C10").Value = "B2C & Sonia"
$ws.Range("C18").Value = "Sonia"

# --- Row 19 (Nov 18): application tag + task description updated ---
$ws.Range("C19").Value = "Sonia"
$ws.Range("D19").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing & Retesting on B2C/B2B app and Sonia Application (Hosyou )"

# --- Row 20 (Nov 19): was blank, now filled in like the rest of the table ---
$ws.Range("C20").Value = "Sonia"
$ws.Range("D20").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing & Retesting on B2C/B2B app and Sonia Application (Soukastu and Hosyou )"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Completed"

# --- Row 21 (Nov 20): was blank, now filled in like the rest of the table ---
$ws.Range("C21").Value = "B2C & B2B"
$ws.Range("D21").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing, Retesting on B2C/B2B app"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "Completed"

# --- View state: scrolled up one row, selection moved to C27 ---
$ws.Range("C27").Select()
